$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 185 (weekly update adds one new record), shifting
# existing rows 185-192 down to 186-193.
$ws.Rows(185).Insert()

$row = 185
$ws.Cells.Item($row, 1).Value = 6
$ws.Cells.Item($row, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($row, 3).Value = "Metropolitana"
$ws.Cells.Item($row, 4).Value = 44568
$ws.Cells.Item($row, 5).Value = 13
$ws.Cells.Item($row, 6).Value = 100112022
$ws.Cells.Item($row, 7).Value = "Arveja Verde"
$ws.Cells.Item($row, 8).Value = "Perfection"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 220
$ws.Cells.Item($row, 11).Value = 14000
$ws.Cells.Item($row, 12).Value = 15000
$ws.Cells.Item($row, 13).Value = 14545
$ws.Cells.Item($row, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item($row, 15).Value = "Carahue"
$ws.Cells.Item($row, 16).Value = 582
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
